# Edit styles of document and bibliography, add references in english,
# other improvements.
#
# This script reapplies a set of paragraph/character style tweaks to the
# active document's style sheet (Heading 1/2, Title, Subtitle, Author,
# Abstract, Bibliography, Table Caption), matching the target OOXML.
#
# NOTE: spacing / indentation values in the OOXML are expressed in twips
# (twentieths of a point), while the Word object model's ParagraphFormat
# properties (SpaceBefore, SpaceAfter, LeftIndent, RightIndent,
# FirstLineIndent) are expressed in points. We therefore convert
# twips -> points by dividing by 20.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Heading 1 (styleId="1")
#   spacing: before 567 -> 340 twips, after 283 -> 340 twips
# ---------------------------------------------------------------------
$s = $d.Styles("Heading 1")
$s.ParagraphFormat.SpaceBefore = 340 / 20
$s.ParagraphFormat.SpaceAfter = 340 / 20

# ---------------------------------------------------------------------
# Heading 2 (styleId="2")
#   spacing: before 200 -> 113 twips, after stays 0
#   indentation: add left=0 right=0 firstLine=709 twips
#   font: Calibri -> Times New Roman, add italic, color -> auto, sz 28 -> 24 (14pt -> 12pt)
# ---------------------------------------------------------------------
$s = $d.Styles("Heading 2")
$s.ParagraphFormat.SpaceBefore = 113 / 20
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LeftIndent = 0
$s.ParagraphFormat.RightIndent = 0
$s.ParagraphFormat.FirstLineIndent = 709 / 20
$s.Font.Name = "Times New Roman"
$s.Font.Italic = $true
$s.Font.ColorIndex = 0
$s.Font.Size = 12

# ---------------------------------------------------------------------
# Title (styleId="Style15")
#   spacing: after 283 -> 0
#   alignment: center -> left
#   bold: true -> explicit false
#   sz: 28 -> 24 (14pt -> 12pt)
# ---------------------------------------------------------------------
$s = $d.Styles("Title")
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.Alignment = 0
$s.Font.Bold = $false
$s.Font.Size = 12

# ---------------------------------------------------------------------
# Subtitle (styleId="Style16")
#   spacing: before 0 -> 340, after 0 -> 340 twips
#   font: add bold
# ---------------------------------------------------------------------
$s = $d.Styles("Subtitle")
$s.ParagraphFormat.SpaceBefore = 340 / 20
$s.ParagraphFormat.SpaceAfter = 340 / 20
$s.Font.Bold = $true

# ---------------------------------------------------------------------
# Author (styleId="Author")
#   spacing: after 283 -> 0
#   sz: 28 -> 24 (14pt -> 12pt)
# ---------------------------------------------------------------------
$s = $d.Styles("Author")
$s.ParagraphFormat.SpaceAfter = 0
$s.Font.Size = 12

# ---------------------------------------------------------------------
# Abstract (styleId="Abstract")
#   spacing: after 283 -> 0
#   indentation: firstLine 709 -> hanging 0 (left/right stay 0)
#   sz: 28 -> 24 (14pt -> 12pt)
# ---------------------------------------------------------------------
$s = $d.Styles("Abstract")
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LeftIndent = 0
$s.ParagraphFormat.RightIndent = 0
$s.ParagraphFormat.FirstLineIndent = -0.001
$s.Font.Size = 12

# ---------------------------------------------------------------------
# Bibliography (styleId="Bibliography")
#   indentation: add left=709 right=0 hanging=709 twips
#   sz: 28 -> 24 (14pt -> 12pt)
# ---------------------------------------------------------------------
$s = $d.Styles("Bibliography")
$s.ParagraphFormat.LeftIndent = 709 / 20
$s.ParagraphFormat.RightIndent = 0
$s.ParagraphFormat.FirstLineIndent = -709 / 20
$s.Font.Size = 12

# ---------------------------------------------------------------------
# Table Caption (styleId="TableCaption")
#   spacing: add before=119 after=119 twips
#   alignment: add center
# ---------------------------------------------------------------------
$s = $d.Styles("Table Caption")
$s.ParagraphFormat.SpaceBefore = 119 / 20
$s.ParagraphFormat.SpaceAfter = 119 / 20
$s.ParagraphFormat.Alignment = 1

Write-Host "Styles updated"
